# Applies the weekly update: inserts two new daily price rows for Jengibre
# (Primera / Segunda) right after row 13, shifting the previously existing
# rows 14-59 down to rows 16-61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 14 (this pushes old rows
# 14-59 down to 16-61, inheriting formatting - e.g. the date style on
# column D - from the row above).
$ws.Rows("14:15").Insert()

# --- New row 14 (Primera) ---
$ws.Cells.Item(14, 1).Value = 9
$ws.Cells.Item(14, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44473
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 100114007
$ws.Cells.Item(14, 7).Value = "Jengibre"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 1060
$ws.Cells.Item(14, 11).Value = 14000
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = 14500
$ws.Cells.Item(14, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 1115
$ws.Cells.Item(14, 17).Value = 13
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# --- New row 15 (Segunda) ---
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44473
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = 100114007
$ws.Cells.Item(15, 7).Value = "Jengibre"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 430
$ws.Cells.Item(15, 11).Value = 11000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 11500
$ws.Cells.Item(15, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 885
$ws.Cells.Item(15, 17).Value = 13
$ws.Cells.Item(15, 18).Value = "Hortaliza"
